$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 500
$ws.Range("B3").Value = 600
$ws.Range("B4").Value = 300
$ws.Range("B5").Value = 123
$ws.Range("B6").Value = 324
$ws.Range("B8").Value = 242
